# Auto-generated edit script: updates cached market-derived profit values
# across the Masamune_Profits workbook sheets (scheduled runner refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3709.6123
$ws.Range("I137").Value = 1008.6071
$ws.Range("J137").Value = 7310.952
$ws.Range("K137").Value = 3025.8213
$ws.Range("L137").Value = 21932.856
$ws.Range("M137").Value = -475.8212999999996
$ws.Range("N137").Value = -27032.856
$ws.Range("H138").Value = 1778.61
$ws.Range("I138").Value = 858.44684
$ws.Range("J138").Value = 2594.6038
$ws.Range("K138").Value = 2575.34052
$ws.Range("L138").Value = 7783.8114
$ws.Range("M138").Value = 2564.65948
$ws.Range("N138").Value = -18063.8114
$ws.Range("H140").Value = 38306.445
$ws.Range("J140").Value = 38306.445
$ws.Range("L140").Value = 38306.445
$ws.Range("N140").Value = -48666.445

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9999.969999999999
$ws.Range("I32").Value = 9540.679
$ws.Range("K32").Value = 9540.679
$ws.Range("M32").Value = -9253.679
$ws.Range("H61").Value = 1349.925
$ws.Range("I61").Value = 1076.5807
$ws.Range("J61").Value = 2291.4443
$ws.Range("K61").Value = 1076.5807
$ws.Range("L61").Value = 2291.4443
$ws.Range("M61").Value = -864.5807
$ws.Range("N61").Value = -2715.4443
$ws.Range("H74").Value = 1524.6809
$ws.Range("I74").Value = 1361
$ws.Range("J74").Value = 2899.6
$ws.Range("K74").Value = 1361
$ws.Range("L74").Value = 2899.6
$ws.Range("M74").Value = -487
$ws.Range("N74").Value = -4647.6
$ws.Range("H77").Value = 1524.6809
$ws.Range("I77").Value = 1361
$ws.Range("J77").Value = 2899.6
$ws.Range("K77").Value = 6805
$ws.Range("L77").Value = 14498
$ws.Range("M77").Value = -2437
$ws.Range("N77").Value = -23234
$ws.Range("H136").Value = 1349.925
$ws.Range("I136").Value = 1076.5807
$ws.Range("J136").Value = 2291.4443
$ws.Range("K136").Value = 3229.7421
$ws.Range("L136").Value = 6874.3329
$ws.Range("M136").Value = -679.7420999999999
$ws.Range("N136").Value = -11974.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 39648
$ws.Range("J59").Value = 39648
$ws.Range("L59").Value = 39648
$ws.Range("N59").Value = -41342
$ws.Range("H86").Value = 1790.6923
$ws.Range("I86").Value = 1814.9166
$ws.Range("K86").Value = 1814.9166
$ws.Range("M86").Value = -691.9166
$ws.Range("H89").Value = 1790.6923
$ws.Range("I89").Value = 1814.9166
$ws.Range("K89").Value = 9074.583000000001
$ws.Range("M89").Value = -3458.583000000001
$ws.Range("H134").Value = 3602.2307
$ws.Range("I134").Value = 2646.6667
$ws.Range("K134").Value = 7940.000100000001
$ws.Range("M134").Value = -5405.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49887.5
$ws.Range("J20").Value = 49887.5
$ws.Range("L20").Value = 49887.5
$ws.Range("N20").Value = -50359.5
$ws.Range("H30").Value = 49887.5
$ws.Range("J30").Value = 49887.5
$ws.Range("L30").Value = 49887.5
$ws.Range("N30").Value = -50069.5
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("H82").Value = 30293.25
$ws.Range("J82").Value = 39391
$ws.Range("L82").Value = 39391
$ws.Range("N82").Value = -40113
$ws.Range("H85").Value = 30293.25
$ws.Range("J85").Value = 39391
$ws.Range("L85").Value = 39391
$ws.Range("N85").Value = -41887
$ws.Range("H88").Value = 36822
$ws.Range("J88").Value = 36822
$ws.Range("L88").Value = 36822
$ws.Range("N88").Value = -37634
$ws.Range("H91").Value = 36822
$ws.Range("J91").Value = 36822
$ws.Range("L91").Value = 36822
$ws.Range("N91").Value = -39630
$ws.Range("H128").Value = 49887.5
$ws.Range("J128").Value = 49887.5
$ws.Range("L128").Value = 49887.5
$ws.Range("N128").Value = -59847.5
$ws.Range("H141").Value = 8109.8887
$ws.Range("J141").Value = 8109.8887
$ws.Range("L141").Value = 8109.8887
$ws.Range("N141").Value = -18469.8887
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 356.0357
$ws.Range("J12").Value = 478.17648
$ws.Range("L12").Value = 1434.52944
$ws.Range("N12").Value = -1780.52944

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 17409.875
$ws.Range("J123").Value = 17409.875
$ws.Range("L123").Value = 17409.875
$ws.Range("N123").Value = -22309.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8541.817999999999
$ws.Range("I46").Value = 3492
$ws.Range("J46").Value = 12750
$ws.Range("K46").Value = 3492
$ws.Range("L46").Value = 12750
$ws.Range("M46").Value = -3304
$ws.Range("N46").Value = -13126
$ws.Range("H47").Value = 25877.5
$ws.Range("J47").Value = 25877.5
$ws.Range("L47").Value = 25877.5
$ws.Range("N47").Value = -26857.5
$ws.Range("H52").Value = 25877.5
$ws.Range("J52").Value = 25877.5
$ws.Range("L52").Value = 25877.5
$ws.Range("N52").Value = -26343.5
$ws.Range("H61").Value = 2184.9333
$ws.Range("I61").Value = 2235.125
$ws.Range("J61").Value = 1984.1666
$ws.Range("K61").Value = 2235.125
$ws.Range("L61").Value = 1984.1666
$ws.Range("M61").Value = -2033.125
$ws.Range("N61").Value = -2388.1666
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("H113").Value = 2184.9333
$ws.Range("I113").Value = 2235.125
$ws.Range("J113").Value = 1984.1666
$ws.Range("K113").Value = 2235.125
$ws.Range("L113").Value = 1984.1666
$ws.Range("M113").Value = -65.125
$ws.Range("N113").Value = -6324.1666
$ws.Range("H128").Value = 39093
$ws.Range("J128").Value = 39093
$ws.Range("L128").Value = 39093
$ws.Range("N128").Value = -49053
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 24680.4
$ws.Range("J101").Value = 24680.4
$ws.Range("L101").Value = 24680.4
$ws.Range("N101").Value = -31170.4
$ws.Range("H104").Value = 48115.2
$ws.Range("J104").Value = 48115.2
$ws.Range("L104").Value = 48115.2
$ws.Range("N104").Value = -55103.2
$ws.Range("H132").Value = 1578.5
$ws.Range("I132").Value = 1211.975
$ws.Range("J132").Value = 3411.125
$ws.Range("K132").Value = 3635.925
$ws.Range("L132").Value = 10233.375
$ws.Range("M132").Value = -1105.925
$ws.Range("N132").Value = -15293.375
$ws.Range("H138").Value = 38328.57
$ws.Range("J138").Value = 38328.57
$ws.Range("L138").Value = 38328.57
$ws.Range("N138").Value = -48608.57
